$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value() = '[Raef%Fadel%NULL%0, Austin R%Morrison%NULL%0, Amit%Vahia%NULL%0, Zachary R%Smith%NULL%0, Zohra%Chaudhry%NULL%0, Pallavi%Bhargava%NULL%0, Joseph%Miller%NULL%0, Rachel M%Kenney%NULL%0, George%Alangaden%NULL%0, Mayur S%Ramesh%mramesh1@hfhs.org%0, Varidhi%Nauriyal%NULL%0, Jayanth%Lakshmikanth%NULL%0, Asif%Abdul Hamed%NULL%0, Owais%Nadeem%NULL%0, Kristin%Griebe%NULL%0, Joseph M%Johnson%NULL%0, Patrick%Bradley%NULL%0, Junior%Uduman%NULL%0, Sara%Hegab%NULL%0, Jennifer%Swiderek%NULL%0, Amanda%Godfrey%NULL%0, Jeffrey%Jennings%NULL%0, Jayna%Gardner-Gray%NULL%0, Adam%Ackerman%NULL%0, Jonathan%Lezotte%NULL%0, Joseph%Ruhala%NULL%0, Linoj%Samuel%NULL%0, Robert J%Tibbetts%NULL%0, Indira%Brar%NULL%0, John%McKinnon%NULL%0, Geehan%Suleyman%NULL%0, Nicholas%Yared%NULL%0, Erica%Herc%NULL%0, Jonathan%Williams%NULL%0, Odaliz Abreu%Lanfranco%NULL%0, Anne%Chen%NULL%0, Marcus%Zervos%NULL%0, Eric%Scher%NULL%0]'
$ws.Range("I2").Value() = ''
$ws.Range("J2").Value() = 'Oxford University Press'

# Row 3
$ws.Range("E3").Value() = '[Xiaowei%Fang%NULL%0, Qing%Mei%NULL%2, Tianjun%Yang%NULL%2, Lei%Li%NULL%0, Yinzhong%Wang%NULL%2, Fei%Tong%NULL%2, Shike%Geng%NULL%2, Aijun%Pan%NULL%3]'
$ws.Range("I3").Value() = ''
$ws.Range("J3").Value() = 'The British Infection Association. Published by Elsevier Ltd.'

# Row 4
$ws.Range("E4").Value() = '[Xiaowei%Fang%NULL%0, Qing%Mei%NULL%0, Tianjun%Yang%NULL%0, Lei%Li%NULL%0, Yinzhong%Wang%NULL%0, Fei%Tong%NULL%0, Shike%Geng%NULL%0, Aijun%Pan%NULL%0]'
$ws.Range("I4").Value() = ''
$ws.Range("J4").Value() = 'The British Infection Association. Published by Elsevier Ltd.'

# Row 5
$ws.Range("E5").Value() = '[Yun%Ling%NULL%0, Shui-Bao%Xu%NULL%2, Yi-Xiao%Lin%NULL%2, Di%Tian%NULL%2, Zhao-Qin%Zhu%NULL%2, Fa-Hui%Dai%NULL%2, Fan%Wu%NULL%2, Zhi-Gang%Song%NULL%2, Wei%Huang%NULL%2, Jun%Chen%NULL%0, Bi-Jie%Hu%NULL%2, Sheng%Wang%NULL%2, En-Qiang%Mao%NULL%2, Lei%Zhu%NULL%0, Wen-Hong%Zhang%NULL%2, Hong-Zhou%Lu%NULL%2, Peng%Lyu%NULL%0, Peng%Lyu%NULL%0]'
$ws.Range("I5").Value() = ''
$ws.Range("J5").Value() = 'Wolters Kluwer Health'

# Row 6
$ws.Range("E6").Value() = '[Xiaofan%Lu%NULL%0, Taige%Chen%NULL%0, Yang%Wang%NULL%0, Jun%Wang%NULL%0, Fangrong%Yan%f.r.yan@163.com%0]'
$ws.Range("I6").Value() = ''
$ws.Range("J6").Value() = 'BioMed Central'

# Row 7
$ws.Range("E7").Value() = '[Yin%Wang%NULL%0, Weiwei%Jiang%NULL%0, Qi%He%NULL%0, Cheng%Wang%NULL%0, Baoju%Wang%NULL%0, Pan%Zhou%NULL%0, Nianguo%Dong%dongnianguo63@gmail.com%0, Qiaoxia%Tong%2013xh0859@hust.edu.cn%0]'
$ws.Range("I7").Value() = ''
$ws.Range("J7").Value() = 'Nature Publishing Group UK'

# Row 8
$ws.Range("C8").Value() = 'Unknown Title'
$ws.Range("D8").Value() = 'Unknown Abstract'
$ws.Range("E8").Value() = '[]'
$ws.Range("F8").Value() = 'not found'
$ws.Range("G8").Value() = 'N/A'
$ws.Range("H8").NumberFormat() = "@"
$ws.Range("H8").Value() = '1970-01-01'
$ws.Range("H8").ClearFormats()
$ws.Range("I8").Value() = ''

# Row 10
$ws.Range("E10").Value() = '[Xiaobo%Yang%NULL%0, Yuan%Yu%NULL%0, Jiqian%Xu%NULL%0, Huaqing%Shu%NULL%0, Jia''an%Xia%NULL%0, Hong%Liu%NULL%0, Yongran%Wu%NULL%0, Lu%Zhang%NULL%0, Zhui%Yu%NULL%0, Minghao%Fang%NULL%0, Ting%Yu%NULL%0, Yaxin%Wang%NULL%0, Shangwen%Pan%NULL%0, Xiaojing%Zou%NULL%0, Shiying%Yuan%NULL%0, You%Shang%NULL%0]'
$ws.Range("I10").Value() = ''
$ws.Range("J10").Value() = 'Elsevier Ltd.'

# Row 11
$ws.Range("E11").Value() = '[Lei%Zha%NULL%0, Shirong%Li%NULL%0, Lingling%Pan%NULL%0, Boris%Tefsen%NULL%0, Yeshan%Li%NULL%0, Neil%French%NULL%0, Liyun%Chen%NULL%0, Gang%Yang%NULL%0, Elmer V%Villanueva%Villanueva@xjtlu.edu.cn%0]'
$ws.Range("I11").Value() = ''
$ws.Range("J11").Value() = 'John Wiley and Sons Inc.'

# Row 12
$ws.Range("E12").Value() = '[Wei%Zhou%NULL%0, Yisi%Liu%NULL%2, Dongdong%Tian%NULL%2, Cheng%Wang%NULL%0, Sa%Wang%NULL%2, Jing%Cheng%NULL%2, Ming%Hu%huming74@163.com%0, Minghao%Fang%fangmh@tjh.tjmu.edu.cn%0, Yue%Gao%gaoyue@bmi.ac.cn%0]'
$ws.Range("I12").Value() = ''
$ws.Range("J12").Value() = 'Nature Publishing Group UK'

# Row 13
$ws.Range("E13").Value() = '[Wei%Zhou%NULL%0, Yisi%Liu%NULL%0, Dongdong%Tian%NULL%0, Cheng%Wang%NULL%0, Sa%Wang%NULL%0, Jing%Cheng%NULL%0, Ming%Hu%huming74@163.com%0, Minghao%Fang%fangmh@tjh.tjmu.edu.cn%0, Yue%Gao%gaoyue@bmi.ac.cn%0]'
$ws.Range("I13").Value() = ''
$ws.Range("J13").Value() = 'Nature Publishing Group UK'

# Row 14
$ws.Range("E14").Value() = '[Rong-chang%Chen%NULL%0, Xiao-ping%Tang%NULL%1, Shou-yong%Tan%NULL%1, Bi-ling%Liang%NULL%1, Zhuo-yue%Wan%NULL%1, Ji-qian%Fang%NULL%1, Nanshan%Zhong%NULL%0]'
$ws.Range("I14").Value() = ''
$ws.Range("J14").Value() = 'The American College of Chest Physicians. Published by Elsevier Inc.'

# Row 15
$ws.Range("D15").Value() = 'We describes a case of a critically ill patient with myocarditis and severe acute respiratory distress syndrome related to coronavirus disease-2019. This case highlights management strategies, including the use of corticosteroids, an interleukin-6 inhibitor, and an aldose reductase inhibitor, resulting in complete clinical recovery.
 (Level of Difficulty: Intermediate.
)'
$ws.Range("E15").Value() = '[Justin%Coyle%NULL%0, Efehi%Igbinomwanhia%NULL%1, Alejandro%Sanchez-Nadales%alejandro.sanchez2@advocatehealth.com%1, Sorin%Danciu%NULL%1, Chae%Chu%NULL%1, Nishit%Shah%NULL%1]'
$ws.Range("I15").Value() = ''
$ws.Range("J15").Value() = 'Elsevier'

# Row 16
$ws.Range("D16").Value() = 'Clusters of patients with novel coronavirus disease 2019 (COVID-19) have been successively reported globally.
 Studies show clear person-to-person transmission.
 The average incubation period is 2–14 days, and mostly 3–7 days.
 However, in some patients, this period may be longer.
 Here, we report a familial cluster of COVID-19 where a 47-year-old woman with long-term use of glucocorticoids did not develop any symptoms within the 14-day quarantine period but was confirmed with COVID-19 by tested positive of antibody on day 40 after she left Wuhan.
 Almost at the same time, her father and sister were diagnosed with COVID-19. The results suggest that the long-term use of glucocorticoids might cause atypical infections, a long incubation period, and extra transmission of COVID-19.'
$ws.Range("E16").Value() = '[Yuanyuan%Han%NULL%0, Mao%Jiang%NULL%1, Da%Xia%NULL%1, Lichao%He%NULL%1, Xin%Lv%NULL%1, Xiaohua%Liao%NULL%1, Jie%Meng%mengjie@csu.edu.cn%1]'
$ws.Range("I16").Value() = ''
$ws.Range("J16").Value() = 'Elsevier Inc.'

# Row 17
$ws.Range("E17").Value() = '[Kui%Liu%NULL%0, Yuan-Yuan%Fang%NULL%0, Yan%Deng%NULL%0, Wei%Liu%NULL%0, Mei-Fang%Wang%NULL%0, Jing-Ping%Ma%NULL%0, Wei%Xiao%NULL%0, Ying-Nan%Wang%NULL%0, Min-Hua%Zhong%NULL%0, Cheng-Hong%Li%NULL%0, Guang-Cai%Li%NULL%0, Hui-Guo%Liu%NULL%0, Xiu-Yuan%Hao%NULL%0, Xiu-Yuan%Hao%NULL%0, Pei-Fang%Wei%NULL%0]'
$ws.Range("I17").Value() = ''
$ws.Range("J17").Value() = 'Wolters Kluwer Health'
